$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.352.45'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '3.242.98'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.24'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.42'
$ws.Range('E6').Value = '  -1.29%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.593'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.135'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.421'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '3.809.83'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.15'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').Value = '68.335.20'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('D17').Value = '3.255.01'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.81'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.44'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '392.47'
$ws.Range('E20').Value = '  +4.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.64'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.25'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.513'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.190'
$ws.Range('E26').Value = '  +4.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.57'
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.98'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.67'
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.91'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.08'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.14'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.92'
$ws.Range('E37').Value = '  +3.87%  '
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.57'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.14'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.54'
$ws.Range('E41').Value = '  -4.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.28'
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.47'
$ws.Range('E43').Value = '  -6.23%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0686'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '343.76'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('D46').Value = '2.598.76'
$ws.Range('E46').Value = '  -4.07%  '
$ws.Range('E47').Value = '  -3.94%  '
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.28'
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '31.50'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('E51').Value = '  -1.49%  '
